$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered) from the
# existing H1 header cell onto the two new header cells, without
# touching the values we just set.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data values for column I (I0)
$ws.Range("I2").Value = 6
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 6
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 7
$ws.Range("I8").Value = 9

# New data values for column J (IF)
$ws.Range("J2").Value = 9
$ws.Range("J3").Value = 3
$ws.Range("J4").Value = 7
$ws.Range("J5").Value = 8
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 8
$ws.Range("J8").Value = 9
